# Update cosinor analysis results on the active sheet to reflect the
# re-run CircadiPy simulation analysis values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = [double]"22.78000000000012"
$ws.Range("H2").Value = [double]"1.797932023684464e-16"
$ws.Range("K2").Value = [double]"43.01969945747804"
$ws.Range("L2").Value = "[37.79533019855198, 48.2440687164041]"
$ws.Range("P2").Value = "[1.5786581702723472, 1.8302371615508113]"
$ws.Range("S2").Value = [double]"52.04031647426044"
$ws.Range("T2").Value = "[48.820334019891746, 55.260298928629126]"
$ws.Range("W2").Value = [double]"16.60044044044053"
$ws.Range("X2").Value = [double]"16.14438438438447"
$ws.Range("Y2").Value = [double]"17.05649649649659"

# --- Row 3 ---
$ws.Range("E3").Value = [double]"25.75000000000059"
$ws.Range("H3").Value = [double]"1.797932023684464e-16"
$ws.Range("I3").Value = [double]"0.1247074164350164"
$ws.Range("K3").Value = [double]"42.4471994568849"
$ws.Range("L3").Value = "[36.69561523090423, 48.19878368286556]"
$ws.Range("O3").Value = [double]"-0.1383684452031542"
$ws.Range("P3").Value = "[-0.27673689040630833, 0.0]"
$ws.Range("Q3").Value = [double]"0.05000000000004534"
$ws.Range("R3").Value = [double]"0.05000000000004534"
$ws.Range("S3").Value = [double]"54.95389049250176"
$ws.Range("T3").Value = "[51.32877070632026, 58.57901027868325]"
$ws.Range("W3").Value = [double]"0.5670670670670788"
$ws.Range("X3").Value = [double]"-1.665334536937735e-15"
$ws.Range("Y3").Value = [double]"1.134134134134159"
